$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.337.83"

$ws.Range("D3").Value = "1.878.51"
$ws.Range("E3").Value = "  +0.25%  "

$ws.Range("E4").Value = "  +0.19%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7110"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.22%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "242.49"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.33%  "

$ws.Range("E7").Value = "  +0.20%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.08020"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.09%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3161"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.50%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.98"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.54%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08301"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.48%  "

$ws.Range("D12").Value = "1.918.83"
$ws.Range("E12").Value = "  +2.90%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.251"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.33%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "94.44"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.64%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.7153"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.39%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.353"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +4.23%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008534"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.52%  "

$ws.Range("D18").Value = "29.351.36"
$ws.Range("E18").Value = "  -0.02%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "243.91"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.48%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.29"
$ws.Range("D20").Style = "Normal"

$ws.Range("D21").Value = "2.133.51"
$ws.Range("E21").Value = "  +0.43%  "

$ws.Range("E22").Value = "  +0.21%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.800"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.43%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.005"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.48%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1558"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.98%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.079"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.28%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "162.83"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.16%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.55"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.22%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.507"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.26%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.422"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.03%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.326"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.06%  "

$ws.Range("E32").Value = "  -7.73%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05386"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.89%  "

$ws.Range("E34").Value = "  +0.02%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7690"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.88%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.186"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.62%  "

$ws.Range("E37").Value = "  -0.57%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01887"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.59%  "

$ws.Range("D39").Value = "1.260.70"
$ws.Range("E39").Value = "  +3.56%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.754"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.87%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.519"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.39%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "113.19"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.02%  "

$ws.Range("B43").Value = "Aave"
$ws.Range("C43").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "74.48"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.30%  "

$ws.Range("B44").Value = "TrustWalletToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9058"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.15%  "

$ws.Range("E45").Value = "  +8.13%  "

$ws.Range("E46").Value = "  +0.17%  "

$ws.Range("D47").Value = "2.028.47"
$ws.Range("E47").Value = "  +0.42%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5233"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.42%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.800"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.26%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.456"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.50%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4376"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.33%  "
